$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add exercise for session 06 (week 6, row 7) -> Aufgaben column F
$ws.Range("F7").Value = "exercises/e06.html"

# Add prep link for session 07 (week 7, row 8) -> Vorbereitung column D
$ws.Range("D8").Value = "https://stats.ifp.uni-mainz.de/ba-ccs-track/ia-zeroshot.html"

# Update selection to D8
$ws.Range("D8").Select()
